$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 8 ("line7" / "line8"), pushing the
# existing "extr1".."extr8" rows down from 8-15 to 10-17.
$ws.Rows("8:9").Insert()

# Make sure the new rows carry the same look as the rest of the table
# (bold/boxed index column) by copying the formatting from the row above.
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the two new rows.
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $false

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber the index column (A) for the rows that shifted down so it
# stays a contiguous 0-based sequence.
for ($r = 10; $r -le 17; $r++) {
    $ws.Range("A$r").Value = $r - 2
}
